{"js": "// Replace invoice placeholder values with the new sample/test values.\n// Each entry is [oldText, newText]; some old values (e.g. the signer's\n// name) appear more than once in the document, so we replace every match.\nconst replacements = [\n  [\"Tovstukha Eduard\", \"wq qw wq qw\"],\n  [\"Molodiznaa 12/12q\", \"wq 23/32q\"],\n  [\"32233 Chmelnitskiy\", \"1111111111 eqe\"],\n  [\"Ukraine\", \"wdw\"],\n  [\"ed@ed.com\", \"dwd@dede\"],\n  [\"380985351072\", \"+380984343994\"],\n  [\"swswsw\", \"dedeed\"],\n  [\"Amount USD: 381\", \"Amount USD: -610080234\"],\n  [\"___________________18-4-2020\", \"___________________16-4-2020\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace invoice placeholder values with the new sample/test values.\n# The values live in the first column of the single table, one per row\n# (row 8 / \"Amount USD\" has a trailing blank paragraph that must stay\n# untouched, and the signer's name in row 1 repeats in row 11).\n# We edit the first paragraph's Range.Text in each target cell directly\n# (rather than Find/Replace) so the xml:space=\"preserve\" markup on the\n# existing <w:t> runs is preserved.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$rowValues = @{\n    1  = \"wq qw wq qw\"\n    2  = \"wq 23/32q\"\n    3  = \"1111111111 eqe\"\n    4  = \"wdw\"\n    5  = \"dwd@dede\"\n    6  = \"+380984343994\"\n    7  = \"dedeed\"\n    8  = \"Amount USD: -610080234\"\n    9  = \"___________________16-4-2020\"\n    11 = \"wq qw wq qw\"\n}\n\nforeach ($rowIndex in $rowValues.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $paragraph = $cell.Range.Paragraphs(1)\n    $paragraph.Range.Text = $rowValues[$rowIndex]\n}\n"}
